# Berechnungen.xlsx — AD5293 functions + small bugfixes
#
# 1) "MCP42050 Widerstände": update the POT1/POT2 sample (row 12) from
#    75/75 to 2/2 (formula in D12 recalculates automatically).
# 2) "AD5293 Widerstände": the resistor table is reworked -
#    - label in A1 changes from "Rwa:" to "Rab:"
#    - the C/E formulas are rewritten (no longer "1024 minus", but a direct
#      ratio), and the shared formula in D is extended down through row 13
#    - two more rows are appended: a blank spacer row (12) and a new data
#      point (13) for 99322 Ohm, whose D13 is a typed-in literal (1020)
#      rather than a copied-down formula.
# 3) Selection/active-sheet bookkeeping matches what Excel would have
#    persisted: "AD5293 Widerstände" becomes the active/selected tab, and
#    the last selections on the other touched sheets are updated too.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) MCP42050 Widerstände
# ---------------------------------------------------------------------
$wsMcp = $wb.Worksheets.Item("MCP42050 Widerstände")
$wsMcp.Range("B12").Value = 2
$wsMcp.Range("C12").Value = 2

# ---------------------------------------------------------------------
# 2) AD5293 Widerstände
# ---------------------------------------------------------------------
$wsAd = $wb.Worksheets.Item("AD5293 Widerstände")

# Label above the table: "Rwa:" -> "Rab:"
$wsAd.Range("A1").Value = "Rab:"

# New formulas for the existing rows (3-11): C is now a direct ratio,
# D stays a ROUND() of C but its shared range will grow to row 13 below,
# E is now based on D directly instead of "1024 minus D".
$wsAd.Range("C3").Formula = '=(B3/$B$1)*1024'
$wsAd.Range("E3").Formula = '=(D3/1024)*$B$1'
$wsAd.Range("C4:C11").Formula = '=(B4/$B$1)*1024'
$wsAd.Range("E4:E11").Formula = '=(D4/1024)*$B$1'

# Blank spacer row 12 (keeps the body formatting, no values).
$wsAd.Range("C12:E12").Value = $null

# New row 13: 99322 Ohm data point.
$wsAd.Range("B13").Value = 99322
$wsAd.Range("C13").Formula = '=(B13/$B$1)*1024'
$wsAd.Range("D13").Value = 1020
$wsAd.Range("E13").Formula = '=(D13/1024)*$B$1'

# Extend the D column's shared ROUND() formula down through the new rows.
$wsAd.Range("D3:D13").Formula = '=ROUND(C3,0)'
$wsAd.Range("D13").Value = 1020

# ---------------------------------------------------------------------
# 3) Window / selection bookkeeping
# ---------------------------------------------------------------------
$wsMcp.Range("D13").Select()
$wsAd.Activate()
$wsAd.Range("D13").Select()

$wsMap = $wb.Worksheets.Item("AD5293 MAP")
$wsMap.Range("G5").Select()

$wsAd.Activate()
